$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D12", "D13", "D14", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D36", "D37", "D39", "D40", "D41", "D42", "D43", "D44", "D46", "D47", "D48", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "62.962.61"
$ws.Range("E2").Value = "  +1.25%  "
$ws.Range("D3").Value = "3.066.69"
$ws.Range("E3").Value = "  +1.37%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "539.70"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").Value = "137.08"
$ws.Range("E6").Value = "  +3.99%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "3.059.29"
$ws.Range("E8").Value = "  +1.06%  "
$ws.Range("D9").Value = "0.493"
$ws.Range("E9").Value = "  +2.09%  "
$ws.Range("D10").Value = "0.155"
$ws.Range("E10").Value = "  +2.34%  "
$ws.Range("D11").Value = "6.25"
$ws.Range("E11").Value = "  +2.36%  "
$ws.Range("D12").Value = "0.453"
$ws.Range("E12").Value = "  -0.60%  "
$ws.Range("D13").Value = "0.0000222"
$ws.Range("E13").Value = "  +5.38%  "
$ws.Range("D14").Value = "34.39"
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("D15").Value = "3.567.99"
$ws.Range("E15").Value = "  +1.68%  "
$ws.Range("D16").Value = "62.999.52"
$ws.Range("E16").Value = "  +1.49%  "
$ws.Range("E17").Value = "  +2.20%  "
$ws.Range("D18").Value = "3.070.14"
$ws.Range("E18").Value = "  +1.26%  "
$ws.Range("D19").Value = "6.62"
$ws.Range("E19").Value = "  +1.00%  "
$ws.Range("D20").Value = "468.42"
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("D21").Value = "13.53"
$ws.Range("E21").Value = "  +2.36%  "
$ws.Range("D22").Value = "0.695"
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("D23").Value = "7.02"
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("D24").Value = "78.39"
$ws.Range("E24").Value = "  +1.66%  "
$ws.Range("D25").Value = "12.08"
$ws.Range("E25").Value = "  +1.34%  "
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").Value = "2.69"
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("D28").Value = "7.86"
$ws.Range("E28").Value = "  -4.18%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("D30").Value = "26.07"
$ws.Range("E30").Value = "  +0.85%  "
$ws.Range("D31").Value = "1.15"
$ws.Range("E31").Value = "  +5.91%  "
$ws.Range("D32").Value = "1.87"
$ws.Range("E32").Value = "  -1.36%  "
$ws.Range("D33").Value = "58.75"
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("E34").Value = "  -4.31%  "
$ws.Range("D35").Value = "5.43"
$ws.Range("E35").Value = "  +9.75%  "
$ws.Range("D36").Value = "5.95"
$ws.Range("E36").Value = "  +1.65%  "
$ws.Range("D37").Value = "482.79"
$ws.Range("E37").Value = "  -0.72%  "
$ws.Range("D38").Value = "3.248.16"
$ws.Range("E38").Value = "  +5.16%  "
$ws.Range("D39").Value = "0.0397"
$ws.Range("E39").Value = "  +2.31%  "
$ws.Range("D40").Value = "0.0791"
$ws.Range("E40").Value = "  +1.45%  "
$ws.Range("D41").Value = "0.117"
$ws.Range("E41").Value = "  +2.01%  "
$ws.Range("D42").Value = "8.12"
$ws.Range("E42").Value = "  +2.48%  "
$ws.Range("D43").Value = "2.56"
$ws.Range("E43").Value = "  +2.10%  "
$ws.Range("D44").Value = "0.250"
$ws.Range("E44").Value = "  +1.06%  "
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("D46").Value = "25.16"
$ws.Range("E46").Value = "  +3.24%  "
$ws.Range("D47").Value = "122.82"
$ws.Range("E47").Value = "  +5.19%  "
$ws.Range("D48").Value = "2.00"
$ws.Range("E48").Value = "  -0.20%  "
$ws.Range("E49").Value = "  +2.48%  "
$ws.Range("D50").Value = "0.0₃0519"
$ws.Range("E50").Value = "  +5.15%  "
$ws.Range("D51").Value = "2.00"
$ws.Range("E51").Value = "  +1.79%  "
